# Extend the table with a new "2022" column (P), mirroring the style of
# the existing "2021" column (O), and set the new selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, font, borders, fill, alignment) of
# column O rows 3-5 into column P rows 3-5 so the new cells inherit the
# same styles used for the other year columns.
$ws.Range("O3:O5").Copy()
$ws.Range("P3:P5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data for the 2022 column.
$ws.Range("P3").Value = 2022
$ws.Range("P4").Value = 15
$ws.Range("P5").Value = 2130.4

# Match the recorded selection after the edit.
$ws.Range("P6").Select()
